# Update the "想去人数" (interest count) figures in column F across the
# relevant worksheets, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5429
$ws1.Range("F6").Value = 819
$ws1.Range("F7").Value = 18
$ws1.Range("F8").Value = 346

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 42

# Sheet "全部类型" (All Types) - combined listing
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5429
$ws4.Range("F6").Value = 819
$ws4.Range("F7").Value = 18
$ws4.Range("F8").Value = 42
$ws4.Range("F9").Value = 346
